$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 0
    4  = 2
    5  = 1
    6  = 2
    7  = 1
    8  = 2
    9  = 1
    10 = 0
    11 = 2
    12 = 1
    13 = 1
    14 = 1
    15 = 2
    16 = 2
    17 = 2
    18 = 0
    19 = 0
    20 = 3
    21 = 1
    22 = 1
    23 = 0
    24 = 2
    25 = 1
    26 = 1
    27 = 1
    28 = 2
    29 = 2
    30 = 3
    31 = 1
    32 = 1
    33 = 1
    34 = 1
    35 = 0
    36 = 1
    37 = 3
    38 = 0
    39 = 1
    40 = 1
    41 = 0
    42 = 2
    43 = 1
    44 = 1
    45 = 1
    46 = 1
    47 = 2
    48 = 2
    49 = 2
    50 = 0
    51 = 2
    52 = 1
    53 = 0
    54 = 1
    55 = 1
    56 = 1
    57 = 0
    58 = 3
    59 = 1
    60 = 2
    61 = 3
    62 = 2
    63 = 0
    64 = 1
    65 = 1
    67 = 2
    68 = 1
    69 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
